# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
#
# Updates the "Metadata" sheet (Version/Date/Publisher/Jurisdiction) and
# removes the duplicate "Contact" row, and updates the "Elements" sheet's
# root Extension row (Short/Definition) to describe the Longitude extension.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Metadata"
# ---------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Bump version and refresh the publish date.
$meta.Cells.Item(3, 2).Value = "6.0.0"
$meta.Cells.Item(8, 2).Value = "2022-01-21T20:46:54+00:00"

# Remove the duplicate "Contact / No display for ContactDetail" row
# (the second of the two identical rows, currently row 11).
$meta.Rows.Item(11).Delete()

# Fill in the Publisher value (was blank).
$meta.Cells.Item(9, 2).Value = "Alvearie Team"

# The remaining old "Contact" row (now row 10) becomes the new
# "Jurisdiction" property.
$meta.Cells.Item(10, 1).Value = "Jurisdiction"
$meta.Cells.Item(10, 2).Value = "United States of America"

# ---------------------------------------------------------------
# Sheet "Elements"
# ---------------------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root "Extension" element row: give it a Longitude-specific
# Short/Definition instead of the generic placeholder text.
$elements.Cells.Item(2, 11).Value = "Longitude"
$elements.Cells.Item(2, 12).Value = "Longitude for the address"
